# "this version has both ref and QA env"
# Adds a new "Json Body" header column (F) to the sheet, alongside the
# widened Email-id/url columns used to show the ref + QA environment
# request bodies, and moves the active selection to E2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the JSON request body column.
$ws.Range("F1").Value = "Json Body"

# Give the new Email-id / Json Body columns (E, F) explicit widths.
$ws.Columns.Item(5).ColumnWidth = 26.5
$ws.Columns.Item(6).ColumnWidth = 24.5

# Move the active cell/selection to E2.
$ws.Range("E2").Select() | Out-Null
